$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Save" header in H1, copying the formatting from the adjacent G1 header cell
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Fill the new column H (rows 2-7) with 0 values
$ws.Range("H2:H7").Value = 0

$excel.CutCopyMode = 0
